$wb = $excel.ActiveWorkbook

# Sheet ALC, row 19
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 5952874
$ws.Range("I19").Value = 28571748
$ws.Range("J19").Value = 539.3158
$ws.Range("K19").Value = 28571748
$ws.Range("L19").Value = 539.3158
$ws.Range("M19").Value = -28571573
$ws.Range("N19").Value = -889.3158

# Sheet ALC, row 33
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 1236.2727
$ws.Range("I33").Value = 1199
$ws.Range("J33").Value = 1267.3334
$ws.Range("K33").Value = 1199
$ws.Range("L33").Value = 1267.3334
$ws.Range("M33").Value = -970
$ws.Range("N33").Value = -1725.3334

# Sheet ALC, row 38
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H38").Value = 1211.3529
$ws.Range("I38").Value = 214.75
$ws.Range("J38").Value = 2097.2222
$ws.Range("K38").Value = 644.25
$ws.Range("L38").Value = 6291.6666
$ws.Range("M38").Value = -272.25
$ws.Range("N38").Value = -7035.6666

# Sheet ALC, row 111
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H111").Value = 1501.625
$ws.Range("I111").Value = 1274.0526
$ws.Range("J111").Value = 2366.4
$ws.Range("K111").Value = 3822.1578
$ws.Range("L111").Value = 7099.200000000001
$ws.Range("M111").Value = -755.1578
$ws.Range("N111").Value = -13233.2

# Sheet ALC, row 112
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H112").Value = 1234.4592
$ws.Range("J112").Value = 1300.8462
$ws.Range("L112").Value = 3902.5386
$ws.Range("N112").Value = -6118.5386

# Sheet ALC, row 137
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 2229721
$ws.Range("I137").Value = 2711059.5
$ws.Range("K137").Value = 8133178.5
$ws.Range("M137").Value = -8130628.5

# Sheet ALC, row 138
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 2416.5244
$ws.Range("I138").Value = 1185.1063
$ws.Range("J138").Value = 4070.1428
$ws.Range("K138").Value = 3555.3189
$ws.Range("L138").Value = 12210.4284
$ws.Range("M138").Value = 1584.6811
$ws.Range("N138").Value = -22490.4284

# Sheet ARM, row 122
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 2716.7058
$ws.Range("I122").Value = 1972.2858
$ws.Range("K122").Value = 5916.857400000001
$ws.Range("M122").Value = -3466.857400000001

# Sheet ARM, row 132
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 24392900
$ws.Range("I132").Value = 29413856
$ws.Range("J132").Value = 5401.7144
$ws.Range("K132").Value = 88241568
$ws.Range("L132").Value = 16205.1432
$ws.Range("M132").Value = -88239038
$ws.Range("N132").Value = -21265.1432

# Sheet BSM, row 7
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H7").Value = 10032.125
$ws.Range("I7").Value = 1450.6
$ws.Range("J7").Value = 24334.666
$ws.Range("K7").Value = 1450.6
$ws.Range("L7").Value = 24334.666
$ws.Range("M7").Value = -1337.6
$ws.Range("N7").Value = -24560.666

# Sheet BSM, row 107
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 2165.4583
$ws.Range("I107").Value = 1588
$ws.Range("J107").Value = 2973.9
$ws.Range("K107").Value = 1588
$ws.Range("L107").Value = 2973.9
$ws.Range("M107").Value = 332
$ws.Range("N107").Value = -6813.9

# Sheet CRP, row 5
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H5").Value = 274.36365
$ws.Range("I5").Value = 284
$ws.Range("J5").Value = 266.33334
$ws.Range("K5").Value = 284
$ws.Range("L5").Value = 266.33334
$ws.Range("M5").Value = -172
$ws.Range("N5").Value = -490.33334

# Sheet CRP, row 8
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H8").Value = 6330
$ws.Range("J8").Value = 6330
$ws.Range("L8").Value = 6330
$ws.Range("N8").Value = -6610

# Sheet CRP, row 10
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H10").Value = 4932
$ws.Range("J10").Value = 14628
$ws.Range("L10").Value = 14628
$ws.Range("N10").Value = -14906

# Sheet CRP, row 11
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H11").Value = 34402.2
$ws.Range("J11").Value = 34402.2
$ws.Range("L11").Value = 34402.2
$ws.Range("N11").Value = -34682.2

# Sheet CRP, row 12
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H12").Value = 8017622
$ws.Range("I12").Value = 13333702
$ws.Range("K12").Value = 13333702
$ws.Range("M12").Value = -13333532

# Sheet CRP, row 14
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H14").Value = 18627.75
$ws.Range("J14").Value = 24670.334
$ws.Range("L14").Value = 24670.334
$ws.Range("N14").Value = -25010.334

# Sheet CRP, row 15
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H15").Value = 3000
$ws.Range("J15").Value = 3000
$ws.Range("L15").Value = 3000
$ws.Range("N15").Value = -3340

# Sheet CRP, row 31
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2859197.2
$ws.Range("I31").Value = 3573346.5
$ws.Range("J31").Value = 2600
$ws.Range("K31").Value = 3573346.5
$ws.Range("L31").Value = 2600
$ws.Range("M31").Value = -3573051.5
$ws.Range("N31").Value = -3190

# Sheet CRP, row 34
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 2859197.2
$ws.Range("I34").Value = 3573346.5
$ws.Range("J34").Value = 2600
$ws.Range("K34").Value = 3573346.5
$ws.Range("L34").Value = 2600
$ws.Range("M34").Value = -3573144.5
$ws.Range("N34").Value = -3004

# Sheet CRP, row 63
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H63").Value = 40000
$ws.Range("J63").Value = 40000
$ws.Range("L63").Value = 40000
$ws.Range("N63").Value = -41372

# Sheet CRP, row 66
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H66").Value = 40000
$ws.Range("J66").Value = 40000
$ws.Range("L66").Value = 120000
$ws.Range("N66").Value = -126864

# Sheet CUL, row 56
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H56").Value = 0
$ws.Range("I56").Value = 0
$ws.Range("K56").Value = 0
$ws.Range("M56").ClearContents()

# Sheet LTW, row 68
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 1461.6
$ws.Range("I68").Value = 1153.9131
$ws.Range("J68").Value = 5000
$ws.Range("K68").Value = 1153.9131
$ws.Range("L68").Value = 5000
$ws.Range("M68").Value = -404.9131
$ws.Range("N68").Value = -6498

# Sheet LTW, row 71
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H71").Value = 1461.6
$ws.Range("I71").Value = 1153.9131
$ws.Range("J71").Value = 5000
$ws.Range("K71").Value = 5769.5655
$ws.Range("L71").Value = 25000
$ws.Range("M71").Value = -2025.5655
$ws.Range("N71").Value = -32488

# Sheet LTW, row 122
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 3595.125
$ws.Range("I122").Value = 2785.7144
$ws.Range("J122").Value = 4224.6665
$ws.Range("K122").Value = 8357.143199999999
$ws.Range("L122").Value = 12673.9995
$ws.Range("M122").Value = -5907.143199999999
$ws.Range("N122").Value = -17573.9995

# Sheet WVR, row 39
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H39").Value = 13792.714
$ws.Range("I39").Value = 1500
$ws.Range("J39").Value = 15841.5
$ws.Range("K39").Value = 1500
$ws.Range("L39").Value = 15841.5
$ws.Range("M39").Value = -1087
$ws.Range("N39").Value = -16667.5

# Sheet WVR, row 122
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 911732.2
$ws.Range("I122").Value = 1252081.8
$ws.Range("J122").Value = 4133.3335
$ws.Range("K122").Value = 3756245.4
$ws.Range("L122").Value = 12400.0005
$ws.Range("M122").Value = -3753795.4
$ws.Range("N122").Value = -17300.0005
